$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.002.99"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "1.829.02"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4609"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3863"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07863"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9584"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").Value = "1.862.68"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.651"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.881"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06749"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009914"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "28.019.60"
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.294"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.085"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").Value = "2.107.56"
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.733"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.972"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9344"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09240"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.288"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.314"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.318"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05863"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.51%  "
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.142"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5569"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.864"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5255"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07019"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.145"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.00%  "
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.318"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.18%  "
